# Auto-generated-assisted Excel COM-interop edit script
# Applies cached-value corrections to the Golem_Profits Leve profitability sheets
# (columns H..N) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR worksheets.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "ALC"; Cell = "H80"; Value = 975 }
    @{ Sheet = "ALC"; Cell = "I80"; Value = 750 }
    @{ Sheet = "ALC"; Cell = "K80"; Value = 2250 }
    @{ Sheet = "ALC"; Cell = "M80"; Value = -1252 }
    @{ Sheet = "ALC"; Cell = "H83"; Value = 975 }
    @{ Sheet = "ALC"; Cell = "I83"; Value = 750 }
    @{ Sheet = "ALC"; Cell = "K83"; Value = 6750 }
    @{ Sheet = "ALC"; Cell = "M83"; Value = -1758 }
    @{ Sheet = "ARM"; Cell = "H35"; Value = 1645.75 }
    @{ Sheet = "ARM"; Cell = "I35"; Value = 1645.75 }
    @{ Sheet = "ARM"; Cell = "K35"; Value = 1645.75 }
    @{ Sheet = "ARM"; Cell = "M35"; Value = -1239.75 }
    @{ Sheet = "ARM"; Cell = "H88"; Value = 2724.3333 }
    @{ Sheet = "ARM"; Cell = "I88"; Value = 1575 }
    @{ Sheet = "ARM"; Cell = "J88"; Value = 3299 }
    @{ Sheet = "ARM"; Cell = "K88"; Value = 1575 }
    @{ Sheet = "ARM"; Cell = "L88"; Value = 3299 }
    @{ Sheet = "ARM"; Cell = "M88"; Value = -1169 }
    @{ Sheet = "ARM"; Cell = "N88"; Value = -4111 }
    @{ Sheet = "ARM"; Cell = "H91"; Value = 2724.3333 }
    @{ Sheet = "ARM"; Cell = "I91"; Value = 1575 }
    @{ Sheet = "ARM"; Cell = "J91"; Value = 3299 }
    @{ Sheet = "ARM"; Cell = "K91"; Value = 1575 }
    @{ Sheet = "ARM"; Cell = "L91"; Value = 3299 }
    @{ Sheet = "ARM"; Cell = "M91"; Value = -171 }
    @{ Sheet = "ARM"; Cell = "N91"; Value = -6107 }
    @{ Sheet = "ARM"; Cell = "H110"; Value = 831.5 }
    @{ Sheet = "ARM"; Cell = "I110"; Value = 831.5 }
    @{ Sheet = "ARM"; Cell = "K110"; Value = 831.5 }
    @{ Sheet = "ARM"; Cell = "M110"; Value = 1213.5 }
    @{ Sheet = "BSM"; Cell = "H86"; Value = 3583 }
    @{ Sheet = "BSM"; Cell = "I86"; Value = 0 }
    @{ Sheet = "BSM"; Cell = "J86"; Value = 3583 }
    @{ Sheet = "BSM"; Cell = "K86"; Value = 0 }
    @{ Sheet = "BSM"; Cell = "L86"; Value = 3583 }
    @{ Sheet = "BSM"; Cell = "M86"; Value = $null }
    @{ Sheet = "BSM"; Cell = "N86"; Value = -5829 }
    @{ Sheet = "BSM"; Cell = "H89"; Value = 3583 }
    @{ Sheet = "BSM"; Cell = "I89"; Value = 0 }
    @{ Sheet = "BSM"; Cell = "J89"; Value = 3583 }
    @{ Sheet = "BSM"; Cell = "K89"; Value = 0 }
    @{ Sheet = "BSM"; Cell = "L89"; Value = 17915 }
    @{ Sheet = "BSM"; Cell = "M89"; Value = $null }
    @{ Sheet = "BSM"; Cell = "N89"; Value = -29147 }
    @{ Sheet = "BSM"; Cell = "H107"; Value = 1845.4445 }
    @{ Sheet = "BSM"; Cell = "I107"; Value = 1845.4445 }
    @{ Sheet = "BSM"; Cell = "J107"; Value = 0 }
    @{ Sheet = "BSM"; Cell = "K107"; Value = 1845.4445 }
    @{ Sheet = "BSM"; Cell = "L107"; Value = 0 }
    @{ Sheet = "BSM"; Cell = "M107"; Value = 74.55549999999994 }
    @{ Sheet = "BSM"; Cell = "N107"; Value = $null }
    @{ Sheet = "CRP"; Cell = "H8"; Value = 12750 }
    @{ Sheet = "CRP"; Cell = "I8"; Value = 0 }
    @{ Sheet = "CRP"; Cell = "J8"; Value = 12750 }
    @{ Sheet = "CRP"; Cell = "K8"; Value = 0 }
    @{ Sheet = "CRP"; Cell = "L8"; Value = 12750 }
    @{ Sheet = "CRP"; Cell = "M8"; Value = $null }
    @{ Sheet = "CRP"; Cell = "N8"; Value = -13030 }
    @{ Sheet = "CRP"; Cell = "H13"; Value = 750 }
    @{ Sheet = "CRP"; Cell = "I13"; Value = 500 }
    @{ Sheet = "CRP"; Cell = "J13"; Value = 1000 }
    @{ Sheet = "CRP"; Cell = "K13"; Value = 500 }
    @{ Sheet = "CRP"; Cell = "L13"; Value = 1000 }
    @{ Sheet = "CRP"; Cell = "M13"; Value = -361 }
    @{ Sheet = "CRP"; Cell = "N13"; Value = -1278 }
    @{ Sheet = "CRP"; Cell = "H47"; Value = 12500 }
    @{ Sheet = "CRP"; Cell = "I47"; Value = 12500 }
    @{ Sheet = "CRP"; Cell = "K47"; Value = 12500 }
    @{ Sheet = "CRP"; Cell = "M47"; Value = -11934 }
    @{ Sheet = "CRP"; Cell = "H107"; Value = 83.333336 }
    @{ Sheet = "CRP"; Cell = "J107"; Value = 50 }
    @{ Sheet = "CRP"; Cell = "L107"; Value = 50 }
    @{ Sheet = "CRP"; Cell = "N107"; Value = -3890 }
    @{ Sheet = "CUL"; Cell = "H4"; Value = 1734.76 }
    @{ Sheet = "CUL"; Cell = "I4"; Value = 1193.4839 }
    @{ Sheet = "CUL"; Cell = "J4"; Value = 2617.8948 }
    @{ Sheet = "CUL"; Cell = "K4"; Value = 3580.4517 }
    @{ Sheet = "CUL"; Cell = "L4"; Value = 7853.6844 }
    @{ Sheet = "CUL"; Cell = "M4"; Value = -3468.4517 }
    @{ Sheet = "CUL"; Cell = "N4"; Value = -8077.6844 }
    @{ Sheet = "CUL"; Cell = "H93"; Value = 750 }
    @{ Sheet = "CUL"; Cell = "J93"; Value = 750 }
    @{ Sheet = "CUL"; Cell = "L93"; Value = 2250 }
    @{ Sheet = "CUL"; Cell = "N93"; Value = -5994 }
    @{ Sheet = "CUL"; Cell = "H113"; Value = 234 }
    @{ Sheet = "CUL"; Cell = "I113"; Value = 350 }
    @{ Sheet = "CUL"; Cell = "J113"; Value = 205 }
    @{ Sheet = "CUL"; Cell = "K113"; Value = 1050 }
    @{ Sheet = "CUL"; Cell = "L113"; Value = 615 }
    @{ Sheet = "CUL"; Cell = "M113"; Value = 1120 }
    @{ Sheet = "CUL"; Cell = "N113"; Value = -4955 }
    @{ Sheet = "CUL"; Cell = "H124"; Value = 428.5 }
    @{ Sheet = "CUL"; Cell = "I124"; Value = 428.5 }
    @{ Sheet = "CUL"; Cell = "K124"; Value = 1285.5 }
    @{ Sheet = "CUL"; Cell = "M124"; Value = 3624.5 }
    @{ Sheet = "GSM"; Cell = "H11"; Value = 1216666.6 }
    @{ Sheet = "GSM"; Cell = "I11"; Value = 1325000 }
    @{ Sheet = "GSM"; Cell = "J11"; Value = 1000000 }
    @{ Sheet = "GSM"; Cell = "K11"; Value = 1325000 }
    @{ Sheet = "GSM"; Cell = "L11"; Value = 1000000 }
    @{ Sheet = "GSM"; Cell = "M11"; Value = -1324861 }
    @{ Sheet = "GSM"; Cell = "N11"; Value = -1000278 }
    @{ Sheet = "GSM"; Cell = "H107"; Value = 1974 }
    @{ Sheet = "GSM"; Cell = "I107"; Value = 1356.6666 }
    @{ Sheet = "GSM"; Cell = "J107"; Value = 2591.3333 }
    @{ Sheet = "GSM"; Cell = "K107"; Value = 1356.6666 }
    @{ Sheet = "GSM"; Cell = "L107"; Value = 2591.3333 }
    @{ Sheet = "GSM"; Cell = "M107"; Value = 563.3334 }
    @{ Sheet = "GSM"; Cell = "N107"; Value = -6431.3333 }
    @{ Sheet = "GSM"; Cell = "H113"; Value = 2881 }
    @{ Sheet = "GSM"; Cell = "I113"; Value = 1668.6666 }
    @{ Sheet = "GSM"; Cell = "J113"; Value = 4699.5 }
    @{ Sheet = "GSM"; Cell = "K113"; Value = 1668.6666 }
    @{ Sheet = "GSM"; Cell = "L113"; Value = 4699.5 }
    @{ Sheet = "GSM"; Cell = "M113"; Value = 501.3334 }
    @{ Sheet = "GSM"; Cell = "N113"; Value = -9039.5 }
    @{ Sheet = "GSM"; Cell = "H132"; Value = 2926.3333 }
    @{ Sheet = "GSM"; Cell = "I132"; Value = 2926.3333 }
    @{ Sheet = "GSM"; Cell = "K132"; Value = 8778.999899999999 }
    @{ Sheet = "GSM"; Cell = "M132"; Value = -6248.999899999999 }
    @{ Sheet = "LTW"; Cell = "H10"; Value = 19000000 }
    @{ Sheet = "LTW"; Cell = "I10"; Value = 19000000 }
    @{ Sheet = "LTW"; Cell = "K10"; Value = 19000000 }
    @{ Sheet = "LTW"; Cell = "M10"; Value = -18999860 }
    @{ Sheet = "LTW"; Cell = "H12"; Value = 1470 }
    @{ Sheet = "LTW"; Cell = "I12"; Value = 1500 }
    @{ Sheet = "LTW"; Cell = "J12"; Value = 1462.5 }
    @{ Sheet = "LTW"; Cell = "K12"; Value = 1500 }
    @{ Sheet = "LTW"; Cell = "L12"; Value = 1462.5 }
    @{ Sheet = "LTW"; Cell = "M12"; Value = -1330 }
    @{ Sheet = "LTW"; Cell = "N12"; Value = -1802.5 }
    @{ Sheet = "LTW"; Cell = "H80"; Value = 27128 }
    @{ Sheet = "LTW"; Cell = "I80"; Value = 0 }
    @{ Sheet = "LTW"; Cell = "J80"; Value = 27128 }
    @{ Sheet = "LTW"; Cell = "K80"; Value = 0 }
    @{ Sheet = "LTW"; Cell = "L80"; Value = 27128 }
    @{ Sheet = "LTW"; Cell = "M80"; Value = $null }
    @{ Sheet = "LTW"; Cell = "N80"; Value = -29374 }
    @{ Sheet = "LTW"; Cell = "H83"; Value = 27128 }
    @{ Sheet = "LTW"; Cell = "I83"; Value = 0 }
    @{ Sheet = "LTW"; Cell = "J83"; Value = 27128 }
    @{ Sheet = "LTW"; Cell = "K83"; Value = 0 }
    @{ Sheet = "LTW"; Cell = "L83"; Value = 81384 }
    @{ Sheet = "LTW"; Cell = "M83"; Value = $null }
    @{ Sheet = "LTW"; Cell = "N83"; Value = -92616 }
    @{ Sheet = "LTW"; Cell = "H93"; Value = 1448.4445 }
    @{ Sheet = "LTW"; Cell = "I93"; Value = 1504.625 }
    @{ Sheet = "LTW"; Cell = "J93"; Value = 999 }
    @{ Sheet = "LTW"; Cell = "K93"; Value = 1504.625 }
    @{ Sheet = "LTW"; Cell = "L93"; Value = 999 }
    @{ Sheet = "LTW"; Cell = "M93"; Value = -256.625 }
    @{ Sheet = "LTW"; Cell = "N93"; Value = -3495 }
    @{ Sheet = "WVR"; Cell = "H107"; Value = 421.3 }
    @{ Sheet = "WVR"; Cell = "I107"; Value = 376.57144 }
    @{ Sheet = "WVR"; Cell = "J107"; Value = 525.6667 }
    @{ Sheet = "WVR"; Cell = "K107"; Value = 1129.71432 }
    @{ Sheet = "WVR"; Cell = "L107"; Value = 1577.0001 }
    @{ Sheet = "WVR"; Cell = "M107"; Value = 790.28568 }
    @{ Sheet = "WVR"; Cell = "N107"; Value = -5417.0001 }
    @{ Sheet = "WVR"; Cell = "H126"; Value = 1878.8 }
    @{ Sheet = "WVR"; Cell = "I126"; Value = 1899.25 }
    @{ Sheet = "WVR"; Cell = "K126"; Value = 5697.75 }
    @{ Sheet = "WVR"; Cell = "M126"; Value = -3227.75 }
    @{ Sheet = "WVR"; Cell = "H138"; Value = 50000 }
    @{ Sheet = "WVR"; Cell = "J138"; Value = 50000 }
    @{ Sheet = "WVR"; Cell = "L138"; Value = 50000 }
    @{ Sheet = "WVR"; Cell = "N138"; Value = -60280 }
)

$touchedSheets = @{}
foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
    $touchedSheets[$u.Sheet] = $true
}

Write-Host "Applied $($updates.Count) cell updates across $($touchedSheets.Keys.Count) worksheets."
